{"js": "// feat: preserve newlines in replace_docx_text\n// \"Apples\" -> \"Apples\\nPears\\nGrapes\" everywhere it occurs in the document.\n// Newlines are written as <w:br/> (manual line break) elements rather than\n// new paragraphs, matching the \"\\n\" -> soft-break behavior described in the\n// commit message (the soft break round-trips back to \"\\n\").\n\nconst body = context.document.body;\n\n// Find every occurrence of \"Apples\" in the document body.\nconst results = body.search(\"Apples\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\n// Office.js represents a manual line break (\"Shift+Enter\", <w:br/>) inside\n// run text as \"\\u000b\" (vertical tab). Using insertText(..., \"Replace\") with\n// embedded \"\\u000b\" characters keeps the inserted text + breaks inside the\n// same run (as separate <w:t>/<w:br/> children), instead of splitting the\n// paragraph the way a literal \"\\n\" would.\nconst replacement = \"Apples\\u000bPears\\u000bGrapes\";\n\nconst count = results.items.length;\nfor (let i = 0; i < count; i++) {\n  results.items[i].insertText(replacement, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# feat: preserve newlines in replace_docx_text\n# \"Apples\" -> \"Apples\\nPears\\nGrapes\" everywhere it occurs in the document.\n# Newlines are written as <w:br/> (manual line break) elements rather than\n# new paragraphs -- this mirrors pressing Shift-Enter in Word, which Word's\n# Find/Replace represents as Chr(11) (vbVerticalTab) in Replacement.Text.\n\n$d = $word.ActiveDocument\n\n$replacement = \"Apples\" + [char]11 + \"Pears\" + [char]11 + \"Grapes\"\n\n$find = $d.Content.Find\n$find.Text = \"Apples\"\n$find.MatchCase = $true\n$find.Replacement.Text = $replacement\n$find.Execute([ref]$find.Text, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$find.Replacement.Text, 2)\n"}
